$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 41666990
$ws.Range("I9").Value = 55555624
$ws.Range("K9").Value = 55555624
$ws.Range("M9").Value = -55555455

$ws.Range("H15").Value = 955.8905999999999
$ws.Range("I15").Value = 955.8905999999999
$ws.Range("K15").Value = 2867.6718
$ws.Range("M15").Value = -2698.6718

$ws.Range("H98").Value = 1342.138
$ws.Range("I98").Value = 1342.138
$ws.Range("K98").Value = 1342.138
$ws.Range("M98").Value = 155.8620000000001

$ws.Range("H112").Value = 313647.5
$ws.Range("J112").Value = 417842.6
$ws.Range("L112").Value = 1253527.8
$ws.Range("N112").Value = -1255743.8

$ws.Range("H115").Value = 701.8461
$ws.Range("I115").Value = 701.8461
$ws.Range("K115").Value = 2105.5383
$ws.Range("M115").Value = -538.5383000000002

$ws.Range("H122").Value = 1342.138
$ws.Range("I122").Value = 1342.138
$ws.Range("K122").Value = 4026.414
$ws.Range("M122").Value = -1576.414

$ws.Range("H136").Value = 64851.555
$ws.Range("J136").Value = 69208
$ws.Range("L136").Value = 69208
$ws.Range("N136").Value = -79408

$ws.Range("H137").Value = 631832.8
$ws.Range("I137").Value = 1167.4
$ws.Range("J137").Value = 807017.7
$ws.Range("K137").Value = 3502.2
$ws.Range("L137").Value = 2421053.1
$ws.Range("M137").Value = -952.2000000000003
$ws.Range("N137").Value = -2426153.1

$ws.Range("H138").Value = 10640024
$ws.Range("I138").Value = 1270.8125
$ws.Range("J138").Value = 12822332
$ws.Range("K138").Value = 3812.4375
$ws.Range("L138").Value = 38466996
$ws.Range("M138").Value = 1327.5625
$ws.Range("N138").Value = -38477276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8276.709000000001
$ws.Range("I32").Value = 5010.7383
$ws.Range("K32").Value = 5010.7383
$ws.Range("M32").Value = -4723.7383

$ws.Range("H74").Value = 3651.4
$ws.Range("I74").Value = 2346.6924
$ws.Range("K74").Value = 2346.6924
$ws.Range("M74").Value = -1472.6924

$ws.Range("H77").Value = 3651.4
$ws.Range("I77").Value = 2346.6924
$ws.Range("K77").Value = 11733.462
$ws.Range("M77").Value = -7365.462

$ws.Range("H110").Value = 1476.5883
$ws.Range("I110").Value = 1446.8
$ws.Range("K110").Value = 1446.8
$ws.Range("M110").Value = 598.2

$ws.Range("H122").Value = 4623.8184
$ws.Range("I122").Value = 4854.067
$ws.Range("K122").Value = 14562.201
$ws.Range("M122").Value = -12112.201

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -748
$ws.Range("N31").ClearContents()

$ws.Range("H134").Value = 2137.12
$ws.Range("I134").Value = 1671.5
$ws.Range("J134").Value = 3999.6
$ws.Range("K134").Value = 5014.5
$ws.Range("L134").Value = 11998.8
$ws.Range("M134").Value = -2479.5
$ws.Range("N134").Value = -17068.8

$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2314.85
$ws.Range("I58").Value = 2125
$ws.Range("J58").Value = 2396.2144
$ws.Range("K58").Value = 2125
$ws.Range("L58").Value = 2396.2144
$ws.Range("M58").Value = -1922
$ws.Range("N58").Value = -2802.2144

$ws.Range("H136").Value = 2314.85
$ws.Range("I136").Value = 2125
$ws.Range("J136").Value = 2396.2144
$ws.Range("K136").Value = 6375
$ws.Range("L136").Value = 7188.6432
$ws.Range("M136").Value = -3825
$ws.Range("N136").Value = -12288.6432

$ws.Range("H141").Value = 252598.5
$ws.Range("I141").Value = 24999.5
$ws.Range("J141").Value = 309498.25
$ws.Range("K141").Value = 24999.5
$ws.Range("L141").Value = 309498.25
$ws.Range("M141").Value = -19819.5
$ws.Range("N141").Value = -319858.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 55561100
$ws.Range("J55").Value = 58829390
$ws.Range("L55").Value = 176488170
$ws.Range("N55").Value = -176488524

$ws.Range("H131").Value = 441799.12
$ws.Range("J131").Value = 627285.3
$ws.Range("L131").Value = 1881855.9
$ws.Range("N131").Value = -1891935.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1259997.5
$ws.Range("J7").Value = 1259997.5
$ws.Range("L7").Value = 1259997.5
$ws.Range("N7").Value = -1260221.5

$ws.Range("H8").Value = 1259997.5
$ws.Range("J8").Value = 1259997.5
$ws.Range("L8").Value = 1259997.5
$ws.Range("N8").Value = -1260275.5

$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 15000
$ws.Range("K18").Value = 15000
$ws.Range("M18").Value = -14707

$ws.Range("H102").Value = 899.55554
$ws.Range("J102").Value = 1000
$ws.Range("L102").Value = 1000
$ws.Range("N102").Value = -4244

$ws.Range("H122").Value = 3362.6785
$ws.Range("J122").Value = 8250
$ws.Range("L122").Value = 24750
$ws.Range("N122").Value = -29650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 17000
$ws.Range("I3").Value = 17000
$ws.Range("K3").Value = 17000
$ws.Range("M3").Value = -16888

$ws.Range("H7").Value = 10231.565
$ws.Range("I7").Value = 11074.096
$ws.Range("K7").Value = 11074.096
$ws.Range("M7").Value = -10962.096

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H15").Value = 17000
$ws.Range("I15").Value = 17000
$ws.Range("K15").Value = 17000
$ws.Range("M15").Value = -16830

$ws.Range("H22").Value = 1750
$ws.Range("J22").Value = 1750
$ws.Range("L22").Value = 1750
$ws.Range("N22").Value = -2340

$ws.Range("H27").Value = 1750
$ws.Range("J27").Value = 1750
$ws.Range("L27").Value = 1750
$ws.Range("N27").Value = -1964

$ws.Range("H40").Value = 1794333.6
$ws.Range("I40").Value = 1909.1818
$ws.Range("J40").Value = 6175815.5
$ws.Range("K40").Value = 1909.1818
$ws.Range("L40").Value = 6175815.5
$ws.Range("M40").Value = -1773.1818
$ws.Range("N40").Value = -6176087.5

$ws.Range("H68").Value = 4156
$ws.Range("I68").Value = 4178.2856
$ws.Range("K68").Value = 4178.2856
$ws.Range("M68").Value = -3429.2856

$ws.Range("H71").Value = 4156
$ws.Range("I71").Value = 4178.2856
$ws.Range("K71").Value = 20891.428
$ws.Range("M71").Value = -17147.428

$ws.Range("H122").Value = 6914760.5
$ws.Range("I122").Value = 23497.715
$ws.Range("J122").Value = 25004324
$ws.Range("K122").Value = 70493.145
$ws.Range("L122").Value = 75012972
$ws.Range("M122").Value = -68043.145
$ws.Range("N122").Value = -75017872

$ws.Range("H126").Value = 10231.565
$ws.Range("I126").Value = 11074.096
$ws.Range("K126").Value = 33222.288
$ws.Range("M126").Value = -30752.288

$ws.Range("H132").Value = 2678.24
$ws.Range("I132").Value = 2062.8462
$ws.Range("K132").Value = 6188.5386
$ws.Range("M132").Value = -3658.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 874.5
$ws.Range("J19").Value = 874.5
$ws.Range("L19").Value = 874.5
$ws.Range("N19").Value = -1222.5

$ws.Range("H25").Value = 2500
$ws.Range("J25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("N25").Value = -3086

$ws.Range("H62").Value = 17250
$ws.Range("J62").Value = 17250
$ws.Range("L62").Value = 17250
$ws.Range("N62").Value = -18498

$ws.Range("H65").Value = 17250
$ws.Range("J65").Value = 17250
$ws.Range("L65").Value = 86250
$ws.Range("N65").Value = -92490

$ws.Range("H132").Value = 1613596.1
$ws.Range("I132").Value = 3390.4546
$ws.Range("J132").Value = 2720612.5
$ws.Range("K132").Value = 10171.3638
$ws.Range("L132").Value = 8161837.5
$ws.Range("M132").Value = -7641.363799999999
$ws.Range("N132").Value = -8166897.5
